# This script reproduces, via Word COM interop, a small proofing/text
# fix-up pass over adresaRaspunsDG.docx:
#
#   1. "R O M A N I A" heading paragraph: mark the trailing "A" run as a
#      grammar-check span (<w:proofErr w:type="gramStart/gramEnd"/>).
#   2. Body paragraph: repair a run split that had cut the word
#      "inregistrata" inside the diacritic "i" (", i" + "nregistrata" ->
#      ", " + "inregistrata"), and drop stray bold formatting that had been
#      applied to the "numarlucrare" placeholder run.
#   3. "Red./Dact./..." footer paragraph: mark "Dact" as a grammar-check
#      span too, which requires splitting the trailing run so "./" carries
#      the <w:proofErr w:type="gramEnd"/> boundary before "S.A.E.S.P./..."
#      continues in its own run.
#
# Each paragraph is located with Find (so the script does not depend on
# fragile absolute paragraph indices), then its whole content is replaced
# in one shot via Range.InsertXML with the corrected run/proofErr markup -
# this keeps every other run attribute (rPr, rsid, etc.) byte-for-byte
# identical to the original, only touching what the diff touches.

$d = $word.ActiveDocument

function Get-ParagraphRangeByText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find anchor text: $searchText"
    }
    return $rng.Paragraphs(1).Range
}

function Set-ParagraphXml($paragraphRange, $newParagraphInnerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' + $newParagraphInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $paragraphRange.InsertXML($pkg)
}

# --- 1. "R O M A N I A" title paragraph -----------------------------------
$romaniaPara = Get-ParagraphRangeByText("R O M")
Set-ParagraphXml $romaniaPara '<w:p w14:paraId="50F19BBF" w14:textId="5FF0E1C6" w:rsidR="00AE2C52" w:rsidRDefault="009D2D76" w:rsidP="00AE2C52"><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:after="46" w:line="264" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:spacing w:val="1"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>R</w:t></w:r><w:r><w:rPr><w:b/><w:spacing w:val="4"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>O</w:t></w:r><w:r><w:rPr><w:b/><w:spacing w:val="4"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>M</w:t></w:r><w:r><w:rPr><w:b/><w:spacing w:val="7"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Â</w:t></w:r><w:r><w:rPr><w:b/><w:spacing w:val="3"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>N</w:t></w:r><w:r><w:rPr><w:b/><w:spacing w:val="3"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>I</w:t></w:r><w:r><w:rPr><w:b/><w:spacing w:val="61"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>A</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'

# --- 2. "..., inregistrata la Inspectoratul ... numarlucrare ..." paragraph
$lucrarePara = Get-ParagraphRangeByText("nregistrată")
Set-ParagraphXml $lucrarePara '<w:p w14:paraId="442062E0" w14:textId="5B264061" w:rsidR="002D2DC6" w:rsidRDefault="002D2DC6" w:rsidP="002D2DC6"><w:pPr><w:pStyle w:val="BodyText"/><w:spacing w:line="592" w:lineRule="auto"/><w:ind w:left="1075" w:right="465" w:firstLine="568"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">La </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>solicitarea</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dvs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. din </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00EF1184" w:rsidRPr="0004069B"><w:t>datalucrare</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>înregistrată</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Inspectoratul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Poliţie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Județean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:spacing w:val="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Iaşi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:rPr><w:spacing w:val="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Serviciul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Arme</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Explozivi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>şi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Substanţe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>periculoase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cu nr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00EF1184" w:rsidRPr="00096415"><w:t>numarlucrare</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>în</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>conformitate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cu </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prevederile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00014362"><w:rPr><w:bCs/></w:rPr><w:t>cadrullegalgd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:rPr><w:spacing w:val="-1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>AVIZĂM</w:t></w:r><w:r><w:rPr><w:spacing w:val="58"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00014362"><w:rPr><w:bCs/></w:rPr><w:t>motivsolicitarepj</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:spacing w:val="59"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">pe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>numitul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:spacing w:val="59"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00014362" w:rsidRPr="004E6707"><w:rPr><w:bCs/></w:rPr><w:t>nume01 nume02</w:t></w:r><w:r w:rsidR="00014362"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00014362" w:rsidRPr="00785539"><w:rPr><w:bCs/></w:rPr><w:t>domiciliat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00014362" w:rsidRPr="00785539"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00014362" w:rsidRPr="00785539"><w:rPr><w:bCs/></w:rPr><w:t>în</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00014362" w:rsidRPr="00785539"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00014362" w:rsidRPr="0004069B"><w:t>adresadomiciliu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00014362" w:rsidRPr="00785539"><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">, C.N.P. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00014362" w:rsidRPr="0004069B"><w:t>cnpsolicitant</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>'

# --- 3. "Red./Dact./S.A.E.S.P./DN./DN./2ex." footer paragraph -------------
$redDactPara = Get-ParagraphRangeByText("Red./")
Set-ParagraphXml $redDactPara '<w:p w14:paraId="748E7EC2" w14:textId="77777777" w:rsidR="00732401" w:rsidRDefault="009D2D76" w:rsidP="004B4B1B"><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:before="29"/><w:ind w:left="200"/><w:rPr><w:rFonts w:ascii="Cambria"/><w:b/><w:sz w:val="12"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cambria"/><w:b/><w:sz w:val="12"/></w:rPr><w:t>Red./</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria"/><w:b/><w:sz w:val="12"/></w:rPr><w:t>Dact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria"/><w:b/><w:sz w:val="12"/></w:rPr><w:t>./</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria"/><w:b/><w:sz w:val="12"/></w:rPr><w:t>S.A.E.S.P./DN./DN./2ex.</w:t></w:r></w:p>'
